# Add the 2025-01-06 09:08 resale-number row (row 10) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (date) and D (week, "01") would otherwise be auto-coerced by
# Excel into a date serial / plain number, losing their original textual
# form (e.g. the leading zero in "01"). Force them to Text first so the
# literal strings are preserved, matching the other rows in the sheet.
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"

$ws.Cells.Item(10, 1).Value = "2025-01-06"
$ws.Cells.Item(10, 2).Value = "09:08:38"
$ws.Cells.Item(10, 3).Value = "Monday"
$ws.Cells.Item(10, 4).Value = "01"
$ws.Cells.Item(10, 5).Value = 127390
$ws.Cells.Item(10, 6).Value = 143665
$ws.Cells.Item(10, 7).Value = 168562
$ws.Cells.Item(10, 8).Value = 158276
$ws.Cells.Item(10, 9).Value = -1
$ws.Cells.Item(10, 10).Value = 142111
$ws.Cells.Item(10, 11).Value = -1
$ws.Cells.Item(10, 12).Value = -1
$ws.Cells.Item(10, 13).Value = 192614
$ws.Cells.Item(10, 14).Value = 114950
$ws.Cells.Item(10, 15).Value = 45469
$ws.Cells.Item(10, 16).Value = 28307
$ws.Cells.Item(10, 17).Value = 63670
$ws.Cells.Item(10, 18).Value = -1
$ws.Cells.Item(10, 19).Value = 47571
$ws.Cells.Item(10, 20).Value = -1
